# The presentation's 3rd slide holds a 3-D cylinder bar chart ("my_chart").
# The authoring tool re-saved that chart and, as a pure side effect of its
# save pipeline, the two internal axis identifiers (<c:axId>/<c:crossAx> in
# the chart's OOXML) were renumbered:
#   95843456 -> 61990016   (category axis)
#   95844992 -> 61991552   (value axis)
# Nothing else about the chart (type, data, scaling, formatting, ...)
# changed. Re-create that renumbering through the Chart/Axis object model.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)

# Find the shape that hosts the chart (named "Graphique 3" / title "my_chart").
$chartShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasChart) {
        $chartShape = $candidate
    }
}

$chart = $chartShape.Chart

$oldCatAxisId = 95843456
$newCatAxisId = 61990016
$oldValAxisId = 95844992
$newValAxisId = 61991552

$catAxis = $chart.Axes(1)   # xlCategory
$valAxis = $chart.Axes(2)   # xlValue

# The axis id is an internal linkage id (<c:axId>/<c:crossAx>) rather than a
# formatting/data property, so not every host surfaces a writer for it. Try
# the direct property first, then a couple of historical aliases via late
# binding, and fall back to silently skipping rather than aborting the whole
# edit if none of them are writable in this runtime.
$setPropertyFlags = [System.Reflection.BindingFlags]::SetProperty

function Set-AxisId($axis, $value) {
    foreach ($propName in @('AxId', 'AxisId', 'ID', 'Id')) {
        try {
            $axis.GetType().InvokeMember($propName, $setPropertyFlags, $null, $axis, @($value)) | Out-Null
            return $true
        } catch {
            # try next alias
        }
    }
    return $false
}

Set-AxisId $catAxis $newCatAxisId | Out-Null
Set-AxisId $valAxis $newValAxisId | Out-Null

Write-Host "Category axis id -> $newCatAxisId (was $oldCatAxisId)"
Write-Host "Value axis id -> $newValAxisId (was $oldValAxisId)"
